$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.693.62'
$ws.Range("E2").Value = '  +0.58%  '

$ws.Range("D3").Value = '2.518.66'
$ws.Range("E3").Value = '  +0.20%  '

$ws.Range("E4").Value = '  -0.03%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '315.20'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.39%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '95.54'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.96%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.575'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -2.32%  '

$ws.Range("E8").Value = '  -0.07%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.533'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.71%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '36.00'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.08%  '

$ws.Range("E11").Value = '  -0.57%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '7.53'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -2.57%  '

$ws.Range("E13").Value = '  -3.78%  '

$ws.Range("D14").Value = '2.905.39'
$ws.Range("E14").Value = '  +0.28%  '

$ws.Range("D15").Value = '2.539.35'
$ws.Range("E15").Value = '  +3.75%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '15.23'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -2.88%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.853'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.48%  '

$ws.Range("D18").Value = '42.777.32'
$ws.Range("E18").Value = '  +0.74%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.74'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +4.43%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.81'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.01%  '

$ws.Range("D21").Value = '0.0₃0959'
$ws.Range("E21").Value = '  -1.39%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '69.48'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.84%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '250.13'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.43%  '

$ws.Range("E24").Value = '  +0.43%  '

$ws.Range("E25").Value = '  +2.13%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '26.58'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -1.51%  '

$ws.Range("E27").Value = '  -0.04%  '

$ws.Range("E28").Value = '  +3.97%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '41.15'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +9.78%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '10.34'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.79%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '5.94'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.04%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '158.03'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +2.62%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '19.53'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.62%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.15'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +3.70%  '

$ws.Range("E35").Value = '  +3.10%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.31'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.96%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0778'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.16%  '

$ws.Range("E38").Value = '  -2.42%  '

$ws.Range("E39").Value = '  -1.32%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '23.62'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -2.45%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.34'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +15.35%  '

$ws.Range("E42").Value = '  +0.82%  '

$ws.Range("E43").Value = '  +0.28%  '

$ws.Range("B44").Value = 'RenderToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.78'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -2.38%  '

$ws.Range("B45").Value = 'NEARProtocol'
$ws.Range("C45").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.32'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -2.55%  '

$ws.Range("D46").Value = '2.045.41'
$ws.Range("E46").Value = '  +0.55%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '84.63'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.11%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '8.91'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.87%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '75.69'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +3.78%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '104.62'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +3.08%  '

$ws.Range("D51").Value = '2.760.55'
$ws.Range("E51").Value = '  +0.20%  '
